# "Add color to the excel" - toggles a set of 0/1 "availability"/"health
# profile" indicator cells on the HealthProfiles and Availability sheets,
# and updates the remembered selection on both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# HealthProfiles sheet
# ---------------------------------------------------------------------
$wsHealth = $wb.Worksheets.Item("HealthProfiles")
$wsHealth.Activate()

$wsHealth.Range("C3").Value = 0
$wsHealth.Range("C4").Value = 0

$wsHealth.Range("I8").Select()

# ---------------------------------------------------------------------
# Availability sheet
# ---------------------------------------------------------------------
$wsAvail = $wb.Worksheets.Item("Availability")
$wsAvail.Activate()

$wsAvail.Range("C2").Value = 0

$wsAvail.Range("B3").Value = 0

$wsAvail.Range("B4").Value = 0

$wsAvail.Range("F5").Value = 0

$wsAvail.Range("E6").Value = 0
$wsAvail.Range("F6").Value = 0

$wsAvail.Range("D7").Value = 0
$wsAvail.Range("E7").Value = 0
$wsAvail.Range("F7").Value = 0

$wsAvail.Range("F9").Value = 0

$wsAvail.Range("E10").Value = 0
$wsAvail.Range("F10").Value = 0

$wsAvail.Range("B11").Value = 0

$wsAvail.Range("B12").Value = 0
$wsAvail.Range("C12").Value = 0

$wsAvail.Range("B13").Value = 0
$wsAvail.Range("C13").Value = 0
$wsAvail.Range("D13").Value = 0

$wsAvail.Range("E14").Value = 0
$wsAvail.Range("F14").Value = 0

$wsAvail.Range("E17").Value = 0
$wsAvail.Range("F17").Value = 0

$wsAvail.Range("B18").Value = 0

$wsAvail.Range("D19").Value = 0
$wsAvail.Range("E19").Value = 0
$wsAvail.Range("F19").Value = 0

$wsAvail.Range("A1:G19").Select()
